# "Generate Report for Handback" — stamp fresh handoff/handback timestamps
# for the e997c3e7-e459-4e68-8caf-ba3a1447e4d5.md row (row 5) across the
# Overview, zh-cn and de-de sheets.
#
# The datetime cells use a custom "yyyy-mm-dd HH:mm:ss" number format; we
# re-assert it before writing so the saved workbook keeps the same visual
# formatting as the other date cells in these columns.
$dateFmt = "yyyy-mm-dd HH:mm:ss"

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for row 5
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").NumberFormat = $dateFmt
$wsOverview.Range("G5").Value = "2016-10-26 07:38:07"

# zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").NumberFormat = $dateFmt
$wsZhCn.Range("H5").Value = "2016-10-26 07:37:55"
$wsZhCn.Range("K5").NumberFormat = $dateFmt
$wsZhCn.Range("K5").Value = "2016-10-26 07:38:33"

# de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").NumberFormat = $dateFmt
$wsDeDe.Range("H5").Value = "2016-10-26 07:38:07"
$wsDeDe.Range("K5").NumberFormat = $dateFmt
$wsDeDe.Range("K5").Value = "2016-10-26 07:38:49"
